$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "20.00", "0.0883")
# must be forced to Text format first, otherwise Excel auto-converts the
# entry to a numeric value and mangles/normalizes the display text
# (loses trailing zeros, introduces floating-point noise, etc.).
$ws.Range("D2").Value = "51.830.72"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "2.810.53"
$ws.Range("E3").Value = "  +1.34%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "355.99"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "111.87"
$ws.Range("E6").Value = "  +2.66%  "
$ws.Range("E7").Value = "  +0.97%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.629"
$ws.Range("E9").Value = "  +8.02%  "
$ws.Range("E10").Value = "  +2.39%  "
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.00"
$ws.Range("E13").Value = "  +3.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.77"
$ws.Range("E14").Value = "  +2.92%  "
$ws.Range("D15").Value = "3.252.79"
$ws.Range("E15").Value = "  +1.39%  "
$ws.Range("D16").Value = "2.811.84"
$ws.Range("E16").Value = "  +1.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.941"
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("D18").Value = "51.838.94"
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.65"
$ws.Range("E19").Value = "  +2.93%  "
$ws.Range("E20").Value = "  +3.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.59"
$ws.Range("E21").Value = "  +3.83%  "
$ws.Range("E22").Value = "  +1.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.44"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "268.71"
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("E25").Value = "  +1.76%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  -0.45%  "
$ws.Range("E28").Value = "  -0.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.80"
$ws.Range("E29").Value = "  +13.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.39"
$ws.Range("E30").Value = "  +2.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.25"
$ws.Range("E31").Value = "  +1.67%  "
$ws.Range("B32").Value = "OKB"
$ws.Range("C32").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "52.48"
$ws.Range("E32").Value = "  +1.97%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.13"
$ws.Range("E33").Value = "  +0.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.65"
$ws.Range("E34").Value = "  +9.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0883"
$ws.Range("E35").Value = "  +5.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0445"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.82"
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("E39").Value = "  +3.20%  "
$ws.Range("E40").Value = "  +0.96%  "
$ws.Range("E41").Value = "  +1.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.52"
$ws.Range("E42").Value = "  -0.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.45"
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.98"
$ws.Range("E44").Value = "  +0.69%  "
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.41"
$ws.Range("E46").Value = "  +4.92%  "
$ws.Range("D47").Value = "2.110.54"
$ws.Range("E47").Value = "  +1.16%  "
$ws.Range("E48").Value = "  +5.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.947"
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("E50").Value = "  -1.63%  "
$ws.Range("E51").Value = "  +6.90%  "
